$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.070.24"
$ws.Range("E2").Value = "  +1.68%  "

$ws.Range("D3").Value = "3.565.10"
$ws.Range("E3").Value = "  +6.95%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'240.11"
$ws.Range("E5").Value = "  +3.53%  "

$ws.Range("D6").Value = "'638.16"
$ws.Range("E6").Value = "  +2.86%  "

$ws.Range("D7").Value = "'1.49"
$ws.Range("E7").Value = "  +7.11%  "

$ws.Range("D8").Value = "'0.403"
$ws.Range("E8").Value = "  +3.58%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("E10").Value = "  +9.32%  "

$ws.Range("D11").Value = "3.565.31"
$ws.Range("E11").Value = "  +6.93%  "

$ws.Range("D12").Value = "'43.56"
$ws.Range("E12").Value = "  +3.18%  "

$ws.Range("E13").Value = "  +4.34%  "

$ws.Range("D14").Value = "'6.42"
$ws.Range("E14").Value = "  +7.26%  "

$ws.Range("D15").Value = "4.244.87"
$ws.Range("E15").Value = "  +7.44%  "

$ws.Range("D16").Value = "96.063.12"
$ws.Range("E16").Value = "  +1.95%  "

$ws.Range("D17").Value = "'0.0000255"
$ws.Range("E17").Value = "  +4.19%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.568.99"
$ws.Range("E18").Value = "  +7.06%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'7.97"
$ws.Range("E19").Value = "  -2.06%  "

$ws.Range("D20").Value = "'12.97"
$ws.Range("E20").Value = "  +17.97%  "

$ws.Range("D21").Value = "'18.04"
$ws.Range("E21").Value = "  +4.23%  "

$ws.Range("D22").Value = "'0.507"
$ws.Range("E22").Value = "  +10.42%  "

$ws.Range("D23").Value = "'516.61"
$ws.Range("E23").Value = "  +4.11%  "

$ws.Range("D24").Value = "'3.45"
$ws.Range("E24").Value = "  -2.05%  "

$ws.Range("E25").Value = "  +9.27%  "

$ws.Range("D26").Value = "'0.0000196"
$ws.Range("E26").Value = "  +6.76%  "

$ws.Range("D27").Value = "'93.24"
$ws.Range("E27").Value = "  +3.55%  "

$ws.Range("D28").Value = "'12.40"
$ws.Range("E28").Value = "  +5.20%  "

$ws.Range("D29").Value = "'3.08"
$ws.Range("E29").Value = "  +16.42%  "

$ws.Range("D30").Value = "'0.145"
$ws.Range("E30").Value = "  +5.33%  "

$ws.Range("D31").Value = "'11.58"
$ws.Range("E31").Value = "  +3.33%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").Value = "'0.184"
$ws.Range("E33").Value = "  +5.29%  "

$ws.Range("D34").Value = "'0.995"
$ws.Range("E34").Value = "  -1.05%  "

$ws.Range("E35").Value = "  +6.26%  "

$ws.Range("D36").Value = "'0.570"
$ws.Range("E36").Value = "  +6.87%  "

$ws.Range("D37").Value = "'586.01"
$ws.Range("E37").Value = "  +9.87%  "

$ws.Range("D38").Value = "'7.75"
$ws.Range("E38").Value = "  +4.80%  "

$ws.Range("D39").Value = "'1.46"
$ws.Range("E39").Value = "  +6.11%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.152"
$ws.Range("E40").Value = "  +3.03%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.935"
$ws.Range("E41").Value = "  +7.05%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").Value = "'1.75"
$ws.Range("E43").Value = "  +3.63%  "

$ws.Range("D44").Value = "'0.0429"
$ws.Range("E44").Value = "  +2.67%  "

$ws.Range("D45").Value = "'23.89"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("D46").Value = "'5.64"
$ws.Range("E46").Value = "  +3.94%  "

$ws.Range("D47").Value = "'3.57"
$ws.Range("E47").Value = "  -5.78%  "

$ws.Range("D48").Value = "'2.19"
$ws.Range("E48").Value = "  +3.66%  "

$ws.Range("D49").Value = "'54.07"
$ws.Range("E49").Value = "  +1.23%  "

$ws.Range("D50").Value = "'8.26"
$ws.Range("E50").Value = "  +2.65%  "

$ws.Range("D51").Value = "'3.14"
$ws.Range("E51").Value = "  +2.59%  "
